$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "64.108.08"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.40%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.528.28"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.00%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "585.54"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.19%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "133.01"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.528.92"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("E12").Value = "  -2.02%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.126.24"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.11%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "27.63"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.22%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.118"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.35%  "
$ws.Range("E16").Value = "  -1.63%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.522.97"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.00%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "64.121.00"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.37%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "9.88"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -3.19%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "14.13"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.43%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "5.62"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.55%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "384.64"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("E23").Value = "  -0.94%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.669.12"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "73.96"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +2.75%  "
$ws.Range("E28").Value = "  -0.34%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.51"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.60%  "
$ws.Range("E30").Value = "  +0.03%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "8.31"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  -1.63%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.537.74"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.38%  "
$ws.Range("E34").Value = "  -0.04%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "23.57"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("E36").Value = "  +0.95%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.37"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("E38").Value = "  -1.37%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "6.91"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.03%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "159.23"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  -0.85%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "26.24"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.11%  "
$ws.Range("E44").Value = "  +0.09%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "41.97"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("E46").Value = "  -4.34%  "
$ws.Range("E47").Value = "  -0.67%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.61"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.71%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.464.58"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +2.20%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "6.85"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.908"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
